# 06/01/25 Commit - added tc TC_TM_028_apply_security_network_Firewall
# Adds a new "Firewall" worksheet, positioned immediately before "Sheet1",
# and populates it with the Add Port / Add Program test data rows.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Add($sheet1)
$ws.Name = "Firewall"

# Header row
$ws.Range("A1").Value = "Template Name"
$ws.Range("B1").Value = "Select Tab"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Port Number"
$ws.Range("E1").Value = "Select Protocol"
$ws.Range("F1").Value = "Program Name"
$ws.Range("G1").Value = "Program Path"

# Row 2 - Add Port (TCP)
$ws.Range("A2").Formula = "=MasterTemplate"
$ws.Range("B2").Value = "Add Port"
$ws.Range("C2").Value = "testname"
$ws.Range("D2").Value = "168.128.1"
$ws.Range("E2").Value = "TCP"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "NA"

# Row 3 - Add Port (UDP)
$ws.Range("A3").Formula = "=MasterTemplate"
$ws.Range("B3").Value = "Add Port"
$ws.Range("C3").Value = "testname"
$ws.Range("D3").Value = "168.128.1"
$ws.Range("E3").Value = "UDP"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "NA"

# Row 4 - Add Program
$ws.Range("A4").Formula = "=MasterTemplate"
$ws.Range("B4").Value = "Add Program"
$ws.Range("C4").Value = "NA"
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "testprogramname"
$ws.Range("G4").Value = "C:\Program Files\Internet Explorer\iexplore.exe"

# Formatting: thin box border around the whole table, yellow header fill,
# matching the look of the workbook's other "Select Tab" data sheets
# (e.g. USBDeviceManager / UserManagement).
$table = $ws.Range("A1:G4")
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2

$header = $ws.Range("A1:G1")
$header.Interior.Color = 65535
$header.Interior.Pattern = 1

$ws.Range("A1:G4").EntireColumn.AutoFit() | Out-Null

$ws.Range("B3").Select()
